$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "21.298.35"
Set-TextValue $ws.Range("E2") "  +4.11%  "
Set-TextValue $ws.Range("D3") "1.544.63"
Set-TextValue $ws.Range("E3") "  +4.97%  "
Set-TextValue $ws.Range("D4") "1.008"
Set-TextValue $ws.Range("E4") "  +0.25%  "
Set-TextValue $ws.Range("D5") "0.9640"
Set-TextValue $ws.Range("E5") "  -1.17%  "
Set-TextValue $ws.Range("D6") "281.74"
Set-TextValue $ws.Range("E6") "  +2.40%  "
Set-TextValue $ws.Range("D7") "0.3633"
Set-TextValue $ws.Range("E7") "  -0.40%  "
Set-TextValue $ws.Range("D8") "0.3210"
Set-TextValue $ws.Range("E8") "  +4.60%  "
Set-TextValue $ws.Range("D9") "40.91"
Set-TextValue $ws.Range("E9") "  +2.92%  "
Set-TextValue $ws.Range("D10") "1.109"
Set-TextValue $ws.Range("E10") "  +5.65%  "
Set-TextValue $ws.Range("D11") "0.06902"
Set-TextValue $ws.Range("E11") "  +4.38%  "
Set-TextValue $ws.Range("E12") "  +0.28%  "
Set-TextValue $ws.Range("D13") "5.716"
Set-TextValue $ws.Range("E13") "  +4.77%  "
Set-TextValue $ws.Range("D14") "18.88"
Set-TextValue $ws.Range("E14") "  +4.98%  "
Set-TextValue $ws.Range("D15") "6.391"
Set-TextValue $ws.Range("E15") "  +3.74%  "
Set-TextValue $ws.Range("D16") "0.00001054"
Set-TextValue $ws.Range("E16") "  +2.60%  "
Set-TextValue $ws.Range("D17") "0.9638"
Set-TextValue $ws.Range("E17") "  -2.08%  "
Set-TextValue $ws.Range("D18") "1.542.48"
Set-TextValue $ws.Range("E18") "  +4.75%  "
Set-TextValue $ws.Range("D19") "0.06130"
Set-TextValue $ws.Range("E19") "  +4.29%  "
Set-TextValue $ws.Range("D20") "72.89"
Set-TextValue $ws.Range("E20") "  +5.05%  "
Set-TextValue $ws.Range("D21") "5.749"
Set-TextValue $ws.Range("E21") "  +5.43%  "
Set-TextValue $ws.Range("D22") "15.22"
Set-TextValue $ws.Range("E22") "  +5.74%  "
Set-TextValue $ws.Range("D23") "11.36"
Set-TextValue $ws.Range("E23") "  +3.85%  "
Set-TextValue $ws.Range("D24") "2.329"
Set-TextValue $ws.Range("E24") "  +3.58%  "
Set-TextValue $ws.Range("D25") "21.377.21"
Set-TextValue $ws.Range("E25") "  +4.24%  "
Set-TextValue $ws.Range("D26") "148.44"
Set-TextValue $ws.Range("E26") "  +4.79%  "
Set-TextValue $ws.Range("D27") "2.249"
Set-TextValue $ws.Range("E27") "  +5.02%  "
Set-TextValue $ws.Range("D28") "17.84"
Set-TextValue $ws.Range("E28") "  +3.39%  "
Set-TextValue $ws.Range("D29") "1.711.01"
Set-TextValue $ws.Range("E29") "  +5.05%  "
Set-TextValue $ws.Range("D30") "118.79"
Set-TextValue $ws.Range("E30") "  +4.59%  "
Set-TextValue $ws.Range("E31") "  +3.98%  "
Set-TextValue $ws.Range("D32") "5.271"
Set-TextValue $ws.Range("E32") "  +5.96%  "
Set-TextValue $ws.Range("D33") "0.8627"
Set-TextValue $ws.Range("E33") "  +7.87%  "
Set-TextValue $ws.Range("D34") "0.08042"
Set-TextValue $ws.Range("E34") "  +2.43%  "
Set-TextValue $ws.Range("D35") "1.517"
Set-TextValue $ws.Range("E35") "  -1.74%  "
Set-TextValue $ws.Range("D36") "4.984"
Set-TextValue $ws.Range("E36") "  +5.11%  "
Set-TextValue $ws.Range("D37") "1.206"
Set-TextValue $ws.Range("E37") "  +4.35%  "
Set-TextValue $ws.Range("D38") "0.05872"
Set-TextValue $ws.Range("E38") "  +1.88%  "
Set-TextValue $ws.Range("D39") "0.02120"
Set-TextValue $ws.Range("E39") "  +3.87%  "
Set-TextValue $ws.Range("D40") "7.824"
Set-TextValue $ws.Range("E40") "  +0.92%  "
Set-TextValue $ws.Range("D41") "10.69"
Set-TextValue $ws.Range("E41") "  +2.12%  "
Set-TextValue $ws.Range("D42") "0.1931"
Set-TextValue $ws.Range("E42") "  +3.17%  "
Set-TextValue $ws.Range("D43") "0.9640"
Set-TextValue $ws.Range("E43") "  -1.13%  "
Set-TextValue $ws.Range("D44") "0.5501"
Set-TextValue $ws.Range("E44") "  +4.02%  "
Set-TextValue $ws.Range("D45") "12.62"
Set-TextValue $ws.Range("E45") "  +5.30%  "
Set-TextValue $ws.Range("D46") "3.575"
Set-TextValue $ws.Range("E46") "  +2.45%  "
Set-TextValue $ws.Range("D47") "0.5472"
Set-TextValue $ws.Range("E47") "  +5.57%  "
Set-TextValue $ws.Range("D48") "122.13"
Set-TextValue $ws.Range("E48") "  +3.83%  "
Set-TextValue $ws.Range("D49") "1.881"
Set-TextValue $ws.Range("E49") "  +6.16%  "
Set-TextValue $ws.Range("D50") "0.06622"
Set-TextValue $ws.Range("E50") "  +2.73%  "
Set-TextValue $ws.Range("D51") "69.73"
Set-TextValue $ws.Range("E51") "  +4.78%  "
